$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style from an existing header cell (e.g. H1) to I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Data values for columns I (I0) and J (IF), rows 2-30
$values = @{
    2  = @(8, 8)
    3  = @(7, 7)
    4  = @(7, 8)
    5  = @(8, 8)
    6  = @(5, 5)
    7  = @(8, 8)
    8  = @(8, 8)
    9  = @(7, 7)
    10 = @(8, 8)
    11 = @(9, 9)
    12 = @(9, 9)
    13 = @(9, 9)
    14 = @(8, 8)
    15 = @(9, 9)
    16 = @(8, 8)
    17 = @(9, 9)
    18 = @(8, 8)
    19 = @(7, 7)
    20 = @(8, 8)
    21 = @(8, 8)
    22 = @(8, 8)
    23 = @(8, 8)
    24 = @(2, 2)
    25 = @(7, 7)
    26 = @(5, 5)
    27 = @(6, 6)
    28 = @(5, 5)
    29 = @(5, 5)
    30 = @(3, 3)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
